# Updated cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '30.123.07'
$ws.Range('E2').Value = '  +5.60%  '
$ws.Range('D3').Value = '1.924.43'
$ws.Range('E3').Value = '  +2.64%  '
$ws.Range('E4').Value = '  -0.97%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +3.51%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.0000'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.84%  '
$ws.Range('E7').Value = '  +1.48%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4008'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +2.76%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.08466'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +1.01%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.91'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +2.58%  '
$ws.Range('E11').Value = '  +1.72%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.27'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +4.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '6.341'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.95%  '
$ws.Range('D14').Value = '1.919.98'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.344'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.41%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.000'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -1.02%  '
$ws.Range('E17').Value = '  +3.42%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001117'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.16%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06777'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +0.92%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '18.04'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +1.81%  '
$ws.Range('E21').Value = '  -0.79%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.065'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +2.25%  '
$ws.Range('D23').Value = '30.115.68'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.21'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.15%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.202'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.51%  '
$ws.Range('D26').Value = '2.141.15'
$ws.Range('E26').Value = '  +2.74%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '160.23'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -0.81%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '20.98'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +1.61%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.470'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  +5.11%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '128.84'
$ws.Range('D30').Style = 'Normal'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.077'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +3.16%  '
$ws.Range('E32').Value = '  +1.26%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '6.073'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +4.86%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.653'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.12%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.02502'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +1.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.06604'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.76%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.2225'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +2.86%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.243'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +4.38%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '9.012'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.80%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.207'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.44%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6546'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +1.79%  '
$ws.Range('E42').Value = '  -0.19%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '11.39'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('E44').Value = '  +1.65%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '13.20'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +1.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.754'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +1.69%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.057'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +2.34%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '125.79'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +3.02%  '
$ws.Range('E49').Value = '  +2.08%  '
$ws.Range('E50').Value = '  +3.65%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.148'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -1.39%  '
